$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.264.31'
$ws.Range("E2").Value = '  +1.72%  '

# Row 3
$ws.Range("D3").Value = '1.893.48'
$ws.Range("E3").Value = '  -1.20%  '

# Row 4
$ws.Range("E4").Value = '  +0.11%  '

# Row 5
$ws.Range("D5").Value = '''324.26'
$ws.Range("E5").Value = '  +2.16%  '

# Row 6
$ws.Range("E6").Value = '  +0.04%  '

# Row 7
$ws.Range("D7").Value = '''0.5180'
$ws.Range("E7").Value = '  -0.30%  '

# Row 8
$ws.Range("D8").Value = '''0.4017'
$ws.Range("E8").Value = '  +1.08%  '

# Row 9
$ws.Range("D9").Value = '''0.08412'
$ws.Range("E9").Value = '  -1.51%  '

# Row 10
$ws.Range("D10").Value = '''42.70'
$ws.Range("E10").Value = '  +0.21%  '

# Row 11
$ws.Range("E11").Value = '  -0.75%  '

# Row 12
$ws.Range("D12").Value = '''23.18'
$ws.Range("E12").Value = '  +11.06%  '

# Row 13
$ws.Range("D13").Value = '''6.439'
$ws.Range("E13").Value = '  +2.09%  '

# Row 14
$ws.Range("D14").Value = '1.899.35'
$ws.Range("E14").Value = '  -0.86%  '

# Row 15
$ws.Range("D15").Value = '''7.317'
$ws.Range("E15").Value = '  -0.56%  '

# Row 16
$ws.Range("D16").Value = '''1.001'
$ws.Range("E16").Value = '  +0.07%  '

# Row 17
$ws.Range("D17").Value = '''94.33'
$ws.Range("E17").Value = '  +0.28%  '

# Row 18
$ws.Range("D18").Value = '''0.00001109'
$ws.Range("E18").Value = '  -0.69%  '

# Row 19
$ws.Range("D19").Value = '''0.06644'
$ws.Range("E19").Value = '  -1.49%  '

# Row 20
$ws.Range("D20").Value = '''18.23'
$ws.Range("E20").Value = '  +1.41%  '

# Row 21
$ws.Range("E21").Value = '  +0.08%  '

# Row 22
$ws.Range("D22").Value = '''5.955'
$ws.Range("E22").Value = '  -1.31%  '

# Row 23
$ws.Range("D23").Value = '30.248.63'
$ws.Range("E23").Value = '  +1.66%  '

# Row 24
$ws.Range("E24").Value = '  +0.62%  '

# Row 25
$ws.Range("E25").Value = '  +0.82%  '

# Row 26
$ws.Range("D26").Value = '2.114.69'
$ws.Range("E26").Value = '  -0.53%  '

# Row 27
$ws.Range("D27").Value = '''21.66'
$ws.Range("E27").Value = '  +3.03%  '

# Row 28
$ws.Range("D28").Value = '''161.83'
$ws.Range("E28").Value = '  +1.83%  '

# Row 29
$ws.Range("D29").Value = '''2.343'
$ws.Range("E29").Value = '  -5.10%  '

# Row 30
$ws.Range("D30").Value = '''129.29'
$ws.Range("E30").Value = '  +0.32%  '

# Row 31
$ws.Range("D31").Value = '''1.090'
$ws.Range("E31").Value = '  +0.73%  '

# Row 32
$ws.Range("D32").Value = '''0.1054'
$ws.Range("E32").Value = '  -0.35%  '

# Row 33
$ws.Range("D33").Value = '''6.100'
$ws.Range("E33").Value = '  -1.82%  '

# Row 34
$ws.Range("D34").Value = '''3.741'
$ws.Range("E34").Value = '  +1.36%  '

# Row 35
$ws.Range("E35").Value = '  -0.50%  '

# Row 36
$ws.Range("D36").Value = '''0.06552'
$ws.Range("E36").Value = '  -1.37%  '

# Row 37
$ws.Range("D37").Value = '''5.340'
$ws.Range("E37").Value = '  +2.21%  '

# Row 38
$ws.Range("D38").Value = '''0.2198'
$ws.Range("E38").Value = '  -0.43%  '

# Row 39
$ws.Range("E39").Value = '  -2.36%  '

# Row 40
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '''8.822'
$ws.Range("E40").Value = '  -3.28%  '

# Row 41
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = '''11.81'
$ws.Range("E41").Value = '  +3.79%  '

# Row 42
$ws.Range("E42").Value = '  -0.82%  '

# Row 43
$ws.Range("D43").Value = '''1.229'
$ws.Range("E43").Value = '  -0.75%  '

# Row 44
$ws.Range("D44").Value = '''0.6084'
$ws.Range("E44").Value = '  -0.85%  '

# Row 45
$ws.Range("D45").Value = '''13.28'
$ws.Range("E45").Value = '  +0.45%  '

# Row 46
$ws.Range("D46").Value = '''3.686'
$ws.Range("E46").Value = '  -0.02%  '

# Row 47
$ws.Range("D47").Value = '''2.055'
$ws.Range("E47").Value = '  -0.63%  '

# Row 48
$ws.Range("E48").Value = '  -0.58%  '

# Row 49
$ws.Range("D49").Value = '''124.66'
$ws.Range("E49").Value = '  -0.19%  '

# Row 50
$ws.Range("D50").Value = '''1.158'
$ws.Range("E50").Value = '  -2.53%  '

# Row 51
$ws.Range("D51").Value = '''79.10'
$ws.Range("E51").Value = '  +0.81%  '
